# Refresh the cryptocurrency price/volume snapshot in the worksheet.
#
# The sheet stores every data cell (Coin, Link, Price, Volume(1h)) as plain
# text, even when the text happens to look like a number (e.g. "240.22" or
# "1.848.08"). Writing such a value straight into Range.Value would make
# Excel auto-convert it to a real number, which would corrupt values like
# "1.848.08" (two dots) and silently change formatting/precision for values
# like "240.22". To avoid that, numeric-looking text is written with a
# leading apostrophe, exactly as if a user typed '240.22 into the cell,
# which forces Excel to keep it as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Address,
        [string]$Text
    )

    $trimmed = $Text.Trim()
    $looksNumeric = ($trimmed.Length -gt 0) -and ($trimmed -match '^[+-]?\d+(\.\d+)?$')

    if ($looksNumeric) {
        $ws.Range($Address).Value = "'" + $Text
    } else {
        $ws.Range($Address).Value = $Text
    }
}

Set-CellText "D2" "29.389.84"
Set-CellText "E2" "  +0.00%  "

Set-CellText "D3" "1.848.08"
Set-CellText "E3" "  -0.02%  "

Set-CellText "E4" "  +0.02%  "

Set-CellText "D5" "240.22"
Set-CellText "E5" "  -0.52%  "

Set-CellText "D6" "0.6304"
Set-CellText "E6" "  -0.65%  "

Set-CellText "D7" "1.0000"
Set-CellText "E7" "  -0.01%  "

Set-CellText "D8" "0.07650"
Set-CellText "E8" "  +1.04%  "

Set-CellText "D9" "0.2935"
Set-CellText "E9" "  -0.93%  "

Set-CellText "D10" "24.60"
Set-CellText "E10" "  -0.75%  "

Set-CellText "D11" "0.07739"
Set-CellText "E11" "  +0.10%  "

Set-CellText "D12" "1.860.09"
Set-CellText "E12" "  -6.29%  "

Set-CellText "D13" "0.00001118"
Set-CellText "E13" "  +12.12%  "

Set-CellText "E14" "  -0.11%  "

Set-CellText "E15" "  -0.58%  "

Set-CellText "D16" "83.70"

Set-CellText "D17" "2.112.48"
Set-CellText "E17" "  -6.71%  "

Set-CellText "D18" "6.191"
Set-CellText "E18" "  +0.89%  "

Set-CellText "D19" "29.408.52"
Set-CellText "E19" "  -0.02%  "

Set-CellText "D20" "229.05"
Set-CellText "E20" "  -0.95%  "

Set-CellText "E21" "  +0.23%  "

Set-CellText "E22" "  +0.04%  "

Set-CellText "D23" "7.516"
Set-CellText "E23" "  -0.53%  "

Set-CellText "E24" "  +0.00%  "

Set-CellText "D25" "157.37"
Set-CellText "E25" "  +0.59%  "

Set-CellText "E26" "  +0.09%  "

Set-CellText "D27" "8.356"
Set-CellText "E27" "  -0.36%  "

Set-CellText "E28" "  -0.38%  "

Set-CellText "D29" "1.466"
Set-CellText "E29" "  -0.42%  "

Set-CellText "D30" "1.299"
Set-CellText "E30" "  +3.55%  "

Set-CellText "D31" "0.05590"
Set-CellText "E31" "  -1.82%  "

Set-CellText "D32" "4.119"
Set-CellText "E32" "  -0.36%  "

Set-CellText "D33" "4.035"
Set-CellText "E33" "  +0.20%  "

Set-CellText "E34" "  +0.28%  "

Set-CellText "D35" "1.158"
Set-CellText "E35" "  +0.08%  "

Set-CellText "D36" "0.7120"
Set-CellText "E36" "  -0.66%  "

Set-CellText "D37" "2.586"
Set-CellText "E37" "  -0.49%  "

Set-CellText "D38" "1.242.02"
Set-CellText "E38" "  -0.16%  "

Set-CellText "D39" "0.01807"
Set-CellText "E39" "  -0.20%  "

Set-CellText "D40" "2.774"
Set-CellText "E40" "  -1.02%  "

Set-CellText "D41" "6.415"
Set-CellText "E41" "  +5.30%  "

Set-CellText "D42" "0.9034"
Set-CellText "E42" "  -0.08%  "

Set-CellText "D43" "0.9998"
Set-CellText "E43" "  +0.00%  "

Set-CellText "D44" "101.93"
Set-CellText "E44" "  +0.05%  "

Set-CellText "D45" "66.09"
Set-CellText "E45" "  -0.18%  "

Set-CellText "E46" "  +1.30%  "

Set-CellText "D47" "7.162"
Set-CellText "E47" "  +1.21%  "

Set-CellText "E48" "  -0.22%  "

# Rows 49 and 50 swap places (EnergySwap now ranks above RenderToken) and
# both get refreshed price/volume figures.
Set-CellText "B49" "EnergySwap"
Set-CellText "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D49" "9.043"
Set-CellText "E49" "  -0.92%  "

Set-CellText "B50" "RenderToken"
Set-CellText "C50" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D50" "1.687"
Set-CellText "E50" "  -1.33%  "

Set-CellText "D51" "0.1122"
Set-CellText "E51" "  -0.42%  "
